# Daily attendance processing - 2025-10-06 15:42:31
#
# The "Recorded By" column (G) on the active sheet lists the users who
# recorded/touched each attendance row, separated by ", ". For every row
# whose list currently starts with "System", the recording order is
# flipped (the list is reversed) so "System" no longer sorts first.
#
# Only the rows listed below (the ones starting with "System,") are
# affected; every other "Recorded By" cell is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,11,12,13,14,15,29,30,32,33,38,39,40,41,42,56,57,58,59,60,65,66,67,68,69,84,85,89,93,110,111,115,119,136,137,141,145)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = [string]$cell.Value2

    if ($current -like "System,*") {
        $parts = $current -split ",\s*"
        $count = $parts.Count

        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $cell.Value = [string]::Join(", ", $reversed)
    }
}
